$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, shifting existing rows 22-78 down to 23-79
$ws.Rows("22:22").Insert()

# Populate the newly inserted row 22 with the new data record
$ws.Cells.Item(22, 1).Value = 9
$ws.Cells.Item(22, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(22, 3).Value = "Metropolitana"
$ws.Cells.Item(22, 4).Value = 44519
$ws.Cells.Item(22, 5).Value = 13
$ws.Cells.Item(22, 6).Value = 100112022
$ws.Cells.Item(22, 7).Value = "Arveja Verde"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 34
$ws.Cells.Item(22, 11).Value = 15000
$ws.Cells.Item(22, 12).Value = 16000
$ws.Cells.Item(22, 13).Value = 15500
$ws.Cells.Item(22, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(22, 15).Value = "Región Metropolitana"
$ws.Cells.Item(22, 16).Value = 620
$ws.Cells.Item(22, 17).Value = 25
$ws.Cells.Item(22, 18).Value = "Hortaliza"
